$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header D1: "Image URL" -> "Image_URL"
$ws.Range("D1").Value = 'Image_URL'

# Update row 2 product data
$ws.Range("A2").Value = '[''Stanley Quencher H2.0 FlowState Stainless Steel Vacuum Insulated Tumbler with Lid and Straw for Water, Iced Tea or Coffee'']'
$ws.Range("B2").Value = '[''4.6'']'
$ws.Range("C2").Value = '[''$45.00'']'
$ws.Range("D2").Value = '[''https://m.media-amazon.com/images/I/41Ewz2W3GPL._AC_SX425_.jpg'']'
$ws.Range("E2").Value = '["About this itemYOUR DREAM TUMBLER: Whichever way your day flows, the H2.0 FlowState tumbler keeps you refreshed with fewer refills. Double wall vacuum insulation means drinks stay cold, iced or hot for hours. Choose between our 14oz, 20oz, 30oz,40oz and 64oz options depending on your hydration needs. The narrow base on all sizes (except 64oz) fits just about any car cup holder, keeping it right by your side.ADVANCED LID CONSTRUCTION: Whether you prefer small sips or maximum thirst quenching, Stanley has developed an advanced FlowState lid, featuring a rotating cover with three positions: a straw opening designed to resist splashes with a seal that holds the reusable straw in place, a drink opening, and a full-cover top for added leak resistance. We’ve also included an ergonomic, comfort-grip handle, so you can easily carry your ice-cold water to work, meetings, the gym or trips out of town.EARTH-FRIENDLY DURABILITY: Constructed of 90% recycled BPA free stainless steel for sustainable sipping, the Stanley Quencher H2.0 has the durability to stand up to a lifetime of use. Eliminate the use of single-use plastic bottles and straws with a travel tumbler built with sustainability in mind.DISHWASHER SAFE: Spend less time hunched over the sink and more time doing the things you love. Cleaning your tumbler and lid couldn''t be easier, just pop them into the dishwasher. Unlike plastic bottles that retain stains & smells, this metallic beauty comes out pristineLIFETIME WARRANTY: Since 1913 we’ve promised to provide rugged, capable gear for food and drink - accessories built to last a lifetime. It’s a promise we still keep. Stanley products purchased from Stanley Resellers come with a lifetime warranty. Rest easy knowing we’ve got your back through it all."]'
$ws.Range("F2").Value = '[''Insulation, Color, Size, Appearance, Coldness, Quality, Value, Condition'']'
